# Update countries & provincias Spain
# Applies the 24-Sep-2020 06:06 data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp in A1
#  - Updates case counters for the countries whose figures changed
#  - Two pairs of rows swap order (their totals crossed one another),
#    so both the country label and the figures for those rows are
#    rewritten to reflect the new ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 06:06"

# --- Helper to write a full data row (A..H) ----------------------------
function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 5 - India (figures updated, same ranking position)
Set-Row 5 "India" 5732518 2334 4674987 966358 0 0 91173

# Rows 21/22 - Pakistan overtakes Turquia
Set-Row 21 "Pakistan" 308217 799 294392 7388 0 5 6437
Set-Row 22 "Turquia" 308069 0 270723 29635 0 0 7711

# Row 37 - Kazajistan (figures updated)
Set-Row 37 "Kazajistan" 107529 79 102323 3507 0 0 1699

# Row 38 - Belgica (figures updated)
Set-Row 38 "Belgica" 106887 1661 19079 77849 0 4 9959

# Rows 158/159 - Belice overtakes Togo
Set-Row 158 "Belice" 1706 37 1019 665 0 1 22
Set-Row 159 "Togo" 1701 0 1297 363 0 0 41

# Row 172 - Islas Turcas y Caicos (figures updated)
Set-Row 172 "Islas Turcas y Caicos" 676 4 588 83 0 0 5

# Row 173 - San Martin (Parte Holandesa) (figures updated)
Set-Row 173 "San Martin (Parte Holandesa)" 616 22 517 78 0 1 21

# Row 188 - Butan (figures updated)
Set-Row 188 "Butan" 261 0 196 65 0 0 0

# Rows 214/215 - Montserrat overtakes Islas Malvinas
Set-Row 214 "Montserrat" 13 0 12 0 0 0 1
Set-Row 215 "Islas Malvinas" 13 0 13 0 0 0 0
